$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the two trailing data rows (Mar-2024 / Apr-2024) that were
# dropped from the series in this revision.
# ------------------------------------------------------------------
$ws.Range("60:61").Delete()

# ------------------------------------------------------------------
# New headers for the added analysis columns.
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Monthly Returns"
$ws.Range("D1").Value = "Average Monthly Returns"
$ws.Range("E1").Value = "Yearly Returns"
$ws.Range("F1").Value = "Average Yearly Returns"

# ------------------------------------------------------------------
# Monthly returns: -B{r}/B{r-1}+1 for every data row (3..59).
# ------------------------------------------------------------------
for ($r = 3; $r -le 59; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=-B$r/B$prev+1"
}

# Average of the monthly returns column.
$ws.Range("D3").Formula = "=AVERAGE(C3:C59)"

# ------------------------------------------------------------------
# Yearly returns, sampled every 12 months.
# ------------------------------------------------------------------
$ws.Range("E21").Formula = "=-B21/B9+1"
$ws.Range("E33").Formula = "=-B33/B21+1"
$ws.Range("E45").Formula = "=-B45/B33+1"
$ws.Range("E57").Formula = "=-B57/B45+1"
$ws.Range("E59").Formula = "=-B59/B47+1"

# Average of the yearly returns column.
$ws.Range("F3").Formula = "=AVERAGE(E21:E59)"

# ------------------------------------------------------------------
# The new formula cells should stay in the "General" format rather
# than inheriting column B's "0.00" number format.
# ------------------------------------------------------------------
$ws.Range("C3:C59").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("E21").ClearFormats()
$ws.Range("E33").ClearFormats()
$ws.Range("E45").ClearFormats()
$ws.Range("E57").ClearFormats()
$ws.Range("E59").ClearFormats()
$ws.Range("F3").ClearFormats()

# ------------------------------------------------------------------
# Restore the selection/cursor position recorded in the saved file.
# ------------------------------------------------------------------
[void]$ws.Range("G45").Select()
